$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 58362
$ws.Range("E2").Value = 981
$ws.Range("F2").Value = 981
$ws.Range("G2").Value = 411
$ws.Range("H2").Value = 633
$ws.Range("I2").Value = 622
$ws.Range("J2").Value = 11
$ws.Range("K2").Value = 82116
$ws.Range("L2").Value = 70925
$ws.Range("M2").Value = 11191
$ws.Range("N2").Value = 10059
$ws.Range("O2").Value = 1133
$ws.Range("P2").Value = 9755
$ws.Range("Q2").Value = 3734
$ws.Range("R2").Value = -3905
$ws.Range("S2").Value = -2179
$ws.Range("T2").Value = 3327
$ws.Range("U2").Value = 407
$ws.Range("V2").Value = 41902
$ws.Range("W2").Value = 1.68
$ws.Range("X2").Value = 1.08
$ws.Range("Y2").Value = 6.53
$ws.Range("Z2").Value = 0.83
$ws.Range("AA2").Value = 633.75
$ws.Range("AB2").Value = -1.82
$ws.Range("AC2").Value = 319
$ws.Range("AD2").Value = 22.37
$ws.Range("AE2").Value = 5156
$ws.Range("AF2").Value = 1.38
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 195101365

# Row 3
$ws.Range("D3").Value = 55407
$ws.Range("E3").Value = 461
$ws.Range("F3").Value = 461
$ws.Range("G3").Value = -1411
$ws.Range("H3").Value = -1392
$ws.Range("I3").Value = -1381
$ws.Range("J3").Value = -11
$ws.Range("K3").Value = 92927
$ws.Range("L3").Value = 84412
$ws.Range("M3").Value = 8516
$ws.Range("N3").Value = 8296
$ws.Range("O3").Value = 219
$ws.Range("P3").Value = 9755
$ws.Range("Q3").Value = 4425
$ws.Range("R3").Value = -8515
$ws.Range("S3").Value = 3289
$ws.Range("T3").Value = 2869
$ws.Range("U3").Value = 1556
$ws.Range("V3").Value = 53476
$ws.Range("W3").Value = 0.83
$ws.Range("X3").Value = -2.51
$ws.Range("Y3").Value = -15.04
$ws.Range("Z3").Value = -1.59
$ws.Range("AA3").Value = 991.23
$ws.Range("AB3").Value = -17
$ws.Range("AC3").Value = -708
$ws.Range("AD3").Value = -6.54
$ws.Range("AE3").Value = 4252
$ws.Range("AF3").Value = 1.09
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 195101365

# Row 4
$ws.Range("D4").Value = 57636
$ws.Range("E4").Value = 2565
$ws.Range("F4").Value = 2565
$ws.Range("G4").Value = 703
$ws.Range("H4").Value = 526
$ws.Range("I4").Value = 493
$ws.Range("J4").Value = 33
$ws.Range("K4").Value = 82292
$ws.Range("L4").Value = 71874
$ws.Range("M4").Value = 10419
$ws.Range("N4").Value = 8920
$ws.Range("O4").Value = 1499
$ws.Range("P4").Value = 10262
$ws.Range("Q4").Value = 6588
$ws.Range("R4").Value = -551
$ws.Range("S4").Value = -5116
$ws.Range("T4").Value = 3749
$ws.Range("U4").Value = 2839
$ws.Range("V4").Value = 46153
$ws.Range("W4").Value = 4.45
$ws.Range("X4").Value = 0.91
$ws.Range("Y4").Value = 5.72
$ws.Range("Z4").Value = 0.6
$ws.Range("AA4").Value = 689.86
$ws.Range("AB4").Value = -11.7
$ws.Range("AC4").Value = 251
$ws.Range("AD4").Value = 16.84
$ws.Range("AE4").Value = 4346
$ws.Range("AF4").Value = 0.97
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 205235294

# Row 5
$ws.Range("D5").Value = 65941
$ws.Range("E5").Value = 2456
$ws.Range("F5").Value = 2456
$ws.Range("G5").Value = 3743
$ws.Range("H5").Value = 2626
$ws.Range("I5").Value = 2519
$ws.Range("J5").Value = -40
$ws.Range("K5").Value = 86565
$ws.Range("L5").Value = 73565
$ws.Range("M5").Value = 12999
$ws.Range("N5").Value = 11458
$ws.Range("O5").Value = 918
$ws.Range("P5").Value = 10262
$ws.Range("Q5").Value = 7244
$ws.Range("R5").Value = -7210
$ws.Range("S5").Value = -157
$ws.Range("T5").Value = 4836
$ws.Range("U5").Value = 2408
$ws.Range("V5").Value = 45180
$ws.Range("W5").Value = 3.72
$ws.Range("X5").Value = 3.98
$ws.Range("Y5").Value = 24.72
$ws.Range("Z5").Value = 3.11
$ws.Range("AA5").Value = 565.91
$ws.Range("AB5").Value = 15.2
$ws.Range("AC5").Value = 1227
$ws.Range("AD5").Value = 3.7
$ws.Range("AE5").Value = 5583
$ws.Range("AF5").Value = 0.8100000000000001
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 205235294

# Row 6
$ws.Range("D6").Value = 71834
$ws.Range("E6").Value = 282
$ws.Range("F6").Value = 282
$ws.Range("G6").Value = -2496
$ws.Range("H6").Value = -1959
$ws.Range("I6").Value = -1979
$ws.Range("K6").Value = 81911
$ws.Range("L6").Value = 70979
$ws.Range("M6").Value = 10932
$ws.Range("N6").Value = 9392
$ws.Range("P6").Value = 10262
$ws.Range("Q6").Value = 7170
$ws.Range("R6").Value = 2231
$ws.Range("S6").Value = -8153
$ws.Range("T6").Value = 3536
$ws.Range("U6").Value = 3634
$ws.Range("V6").Value = 34402
$ws.Range("W6").Value = 0.39
$ws.Range("X6").Value = -2.73
$ws.Range("Y6").Value = -18.98
$ws.Range("Z6").Value = -2.33
$ws.Range("AA6").Value = 649.28
$ws.Range("AB6").Value = -5.98
$ws.Range("AC6").Value = -964
$ws.Range("AD6").Value = -4.29
$ws.Range("AE6").Value = 4576
$ws.Range("AF6").Value = 0.9
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 205235294

# Row 7
$ws.Range("D7").Value = 71000
$ws.Range("E7").Value = -2779
$ws.Range("G7").Value = -7286
$ws.Range("H7").Value = -5725
$ws.Range("I7").Value = -5685
$ws.Range("K7").Value = 107344
$ws.Range("L7").Value = 96104
$ws.Range("M7").Value = 11240
$ws.Range("N7").Value = 9887
$ws.Range("P7").Value = 10660
$ws.Range("Q7").Value = 7193
$ws.Range("R7").Value = -13403
$ws.Range("S7").Value = -1897
$ws.Range("T7").Value = 13243
$ws.Range("U7").Value = 1057
$ws.Range("W7").Value = -3.91
$ws.Range("X7").Value = -8.06
$ws.Range("Y7").Value = -58.97
$ws.Range("Z7").Value = -6.05
$ws.Range("AA7").Value = 854.99
$ws.Range("AC7").Value = -2625
$ws.Range("AD7").Value = -1.74
$ws.Range("AE7").Value = 4429
$ws.Range("AF7").Value = 1.03
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 73220
$ws.Range("E8").Value = -318
$ws.Range("G8").Value = -3665
$ws.Range("H8").Value = -2870
$ws.Range("I8").Value = -2299
$ws.Range("K8").Value = 111908
$ws.Range("L8").Value = 95160
$ws.Range("M8").Value = 16748
$ws.Range("N8").Value = 15435
$ws.Range("P8").Value = 21283
$ws.Range("Q8").Value = 3932
$ws.Range("R8").Value = -3648
$ws.Range("S8").Value = 5703
$ws.Range("T8").Value = 4325
$ws.Range("U8").Value = 1867
$ws.Range("W8").Value = -0.43
$ws.Range("X8").Value = -3.92
$ws.Range("Y8").Value = -18.16
$ws.Range("Z8").Value = -2.62
$ws.Range("AA8").Value = 568.1799999999999
$ws.Range("AC8").Value = -1030
$ws.Range("AD8").Value = -4.43
$ws.Range("AE8").Value = 6914
$ws.Range("AF8").Value = 0.66
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 75502
$ws.Range("E9").Value = 1337
$ws.Range("G9").Value = -1546
$ws.Range("H9").Value = -1214
$ws.Range("I9").Value = -768
$ws.Range("K9").Value = 110896
$ws.Range("L9").Value = 95409
$ws.Range("M9").Value = 15487
$ws.Range("N9").Value = 14216
$ws.Range("P9").Value = 21283
$ws.Range("Q9").Value = 7510
$ws.Range("R9").Value = -3777
$ws.Range("S9").Value = -1224
$ws.Range("T9").Value = 4710
$ws.Range("U9").Value = 1943
$ws.Range("W9").Value = 1.77
$ws.Range("X9").Value = -1.61
$ws.Range("Y9").Value = -5.18
$ws.Range("Z9").Value = -1.09
$ws.Range("AA9").Value = 616.05
$ws.Range("AC9").Value = -344
$ws.Range("AD9").Value = -13.26
$ws.Range("AE9").Value = 6368
$ws.Range("AF9").Value = 0.72
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").ClearContents()
